$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 0.41936842089428694
$ws.Range("A2").Value = -0.0099999990545747153
$ws.Range("A3").Value = -0.0089999990379201478
$ws.Range("A4").Value = -0.011999999761052038
$ws.Range("A5").Value = -0.0059999990508083911
$ws.Range("A6").Value = -0.0059999990240591217
$ws.Range("A7").Value = -0.019999998856523149
$ws.Range("A8").Value = -0.019999998850575906
$ws.Range("A9").Value = -0.0059999990110295442
$ws.Range("A10").Value = -0.0059999990078267729
$ws.Range("A11").Value = -0.0044999990255725208
$ws.Range("A12").Value = -0.0059999990079284693
$ws.Range("A13").Value = -0.0059999990083392518
$ws.Range("A14").Value = -0.011999998938307499
$ws.Range("A15").Value = -0.0059999990111290202
$ws.Range("A16").Value = -0.0059999990149406379
$ws.Range("A17").Value = 0.022899123336109994
$ws.Range("A18").Value = 0.0068734048155345562
$ws.Range("A19").Value = -0.060037646325334482
$ws.Range("A20").Value = -0.0089999990334277413
$ws.Range("A21").Value = -0.063427064677365141
$ws.Range("A22").Value = -0.0089999990265030583
$ws.Range("A23").Value = -0.0089999990292222165
$ws.Range("A24").Value = -0.041999998626818247
$ws.Range("A25").Value = -0.04199999861921544
$ws.Range("A26").Value = -0.005999999023881486
$ws.Range("A27").Value = -0.0059999990202306286
$ws.Range("A28").Value = 0.019159802822729333
$ws.Range("A29").Value = -0.011999998922354038
$ws.Range("A30").Value = -0.019999998821216725
$ws.Range("A31").Value = -0.014999998870619535
$ws.Range("A32").Value = -0.020999998798681752
$ws.Range("A33").Value = -0.0059999989752439475
